$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data row (row 2) and write the new data into row 1.
$ws.Rows(2).ClearContents()

$ws.Range("A1").Value = "value"
$ws.Range("B1").Value = "world!"
